$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header format (from H1, the last existing header cell) onto the two new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header row (row 1): new column titles
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-7 for new columns I and J
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 4

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 2

$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 8

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 3

$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 6
